$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# New row of data to append at row 40
$ws.Cells.Item(40, 1).Value = 41225
$ws.Cells.Item(40, 1).NumberFormat = "ddd\ dd/mm/yyyy"

$ws.Cells.Item(40, 2).Value = 3.5

$ws.Cells.Item(40, 4).Value = "Installer creation scripts continued, new test case tc08"

# Update selection to match the diff (E40)
$ws.Range("E40").Select()
